$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "edited scripts and added job run results"
#
# Row 42 was an empty template row; it now holds a new job-run result that
# failed to produce usable overlap/alpha numbers (hence the "N/A" markers),
# mirroring the layout of the data-row directly above it (row 40).
# ---------------------------------------------------------------------------

$ws.Range("A42").Value = "ukb51139_subset.csv"
$ws.Range("B42").Value = "28012 x 1081"
$ws.Range("C42").Value = "all"
$ws.Range("D42").Value = "no events"
$ws.Range("E42").Value = "> 140/80"
$ws.Range("F42").Value = "zscore"
$ws.Range("G42").Value = "median"
$ws.Range("H42").Value = "none"
$ws.Range("I42").Value = 25
$ws.Range("K42").Value = "N/A"
$ws.Range("L42").Value = "-211.7 & -59.8"
$ws.Range("M42").Value = "47.7 & 45.3"
$ws.Range("N42").Value = "N/A"
$ws.Range("O42").Value = "N/A"

# ---------------------------------------------------------------------------
# The numeric-result columns (I, K, N, O) for the previously-blank rows
# 41-48 were still wearing the old, now-stale "template" number format
# (a leftover theme-coloured font); bring them in line with the style
# already used throughout the rest of the results table (rows 2-40), i.e.
# the explicit-black-font variant. J stays a plain blank separator column.
# ---------------------------------------------------------------------------

$resultRows = 40,41,42,43,44,45,46,47,48
foreach ($r in $resultRows) {
    foreach ($col in @("I", "K", "N")) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "#,##0"
        $cell.HorizontalAlignment = -4152
        $cell.Font.Color = 0
    }
    $oCell = $ws.Range("O$r")
    # Row 42's O cell holds the "N/A" text, so it picked up the same
    # integer-style formatting as the other N/A cells in that row
    # rather than the usual two-decimal numeric style.
    if ($r -eq 42) {
        $oCell.NumberFormat = "#,##0"
    } else {
        $oCell.NumberFormat = "#,##0.00"
    }
    $oCell.HorizontalAlignment = -4152
    $oCell.Font.Color = 0
}

# ---------------------------------------------------------------------------
# Rows 40 & 41 grew to match the taller row height used elsewhere
# (19.5pt instead of 18.75pt).
# ---------------------------------------------------------------------------

$ws.Rows.Item(40).RowHeight = 19.5
$ws.Rows.Item(41).RowHeight = 19.5

# ---------------------------------------------------------------------------
# Columns C, D and I were re-sized (best-fit) by the reporting script for
# the new, narrower content widths.
# ---------------------------------------------------------------------------

$ws.Columns.Item(3).ColumnWidth = 15.333333333333334
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666
$ws.Columns.Item(9).ColumnWidth = 10.833333333333334
